$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Tipo" column header
$ws.Range("C1").Value = "Tipo"

# Existing rows (2-11) are "discreta" type materials
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = "discreta"
}

# New row 12: pelo de vaca - discreta
$ws.Range("A12").Value = "pelo de vaca"
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = "discreta"

# New row 13: pelo de cerdo - continua
$ws.Range("A13").Value = "pelo de cerdo"
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = "continua"
